$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.275898814201355
$ws.Range("B1").Value = 2.303384304046631
$ws.Range("D1").Value = 1.392434120178223
$ws.Range("E1").Value = 0.8453736305236816
